$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.900.36"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.079.62"
$ws.Range("E3").Value = "  -0.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.85"
$ws.Range("E5").Value = "  -1.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.08"
$ws.Range("E6").Value = "  -0.84%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.075.62"
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  +0.83%  "

# Row 10
$ws.Range("E10").Value = "  -0.75%  "

# Row 11
$ws.Range("E11").Value = "  -1.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").Value = "  -0.91%  "

# Row 13
$ws.Range("E13").Value = "  +4.63%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.40"
$ws.Range("E14").Value = "  -1.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.573.56"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.929.87"
$ws.Range("E16").Value = "  +1.15%  "

# Row 17
$ws.Range("E17").Value = "  +0.93%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.080.08"
$ws.Range("E18").Value = "  +0.30%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.66"
$ws.Range("E19").Value = "  +0.37%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.46"
$ws.Range("E20").Value = "  -2.00%  "

# Row 21
$ws.Range("E21").Value = "  -0.76%  "

# Row 22
$ws.Range("E22").Value = "  -0.26%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.09"
$ws.Range("E23").Value = "  -0.69%  "

# Row 24
$ws.Range("E24").Value = "  +2.71%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.14"
$ws.Range("E25").Value = "  -0.85%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.11%  "

# Row 27
$ws.Range("E27").Value = "  -0.87%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.09"
$ws.Range("E28").Value = "  -2.99%  "

# Row 29
$ws.Range("E29").Value = "  +0.10%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.21"
$ws.Range("E30").Value = "  -0.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  +1.58%  "

# Row 32
$ws.Range("E32").Value = "  -2.65%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "57.19"
$ws.Range("E33").Value = "  -4.53%  "

# Row 34
$ws.Range("E34").Value = "  -6.56%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "503.25"
$ws.Range("E35").Value = "  -4.21%  "

# Row 36
$ws.Range("E36").Value = "  +3.77%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.99"
$ws.Range("E37").Value = "  +0.95%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.229.05"
$ws.Range("E38").Value = "  +6.02%  "

# Row 39
$ws.Range("E39").Value = "  -1.38%  "

# Row 40
$ws.Range("E40").Value = "  -0.24%  "

# Row 41
$ws.Range("E41").Value = "  -0.66%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.10"
$ws.Range("E42").Value = "  +0.04%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").Value = "  -1.05%  "

# Row 44
$ws.Range("E44").Value = "  -0.71%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "121.98"
$ws.Range("E46").Value = "  +0.11%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "2.04"
$ws.Range("E47").Value = "  -0.15%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₃0528"
$ws.Range("E48").Value = "  +3.95%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.38"
$ws.Range("E49").Value = "  +0.29%  "

# Row 50
$ws.Range("E50").Value = "  +1.30%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.38"
$ws.Range("E51").Value = "  +3.53%  "
